$d = $word.ActiveDocument

# Update the PHONE contact label to TELEGRAM
$d.Content.Find.Execute("PHONE", $true, $false, $false, $false, $false,
                         $true, 1, $false, "TELEGRAM", 2)

# Update the phone number to the Telegram handle
$d.Content.Find.Execute("067 518 22 22", $true, $false, $false, $false, $false,
                         $true, 1, $false, "@vadymvoitsekhovskyi", 2)

# Best-effort: mark the unused "Default Paragraph Font" style as semi-hidden
# (mirrors Word's automatic "unused built-in style" bookkeeping). Guarded so
# the rest of the script keeps running even where this isn't supported.
try {
    $style = $d.Styles("Default Paragraph Font")
    $style.Hidden = $true
} catch {
}
